# Update countries & provincias Spain
# - Refresh the "last updated" timestamp in A1
# - Refresh COVID stats for Estados Unidos (row 4) and Irlanda (row 34)
# - Refresh COVID stats for Gabon / Croacia; Gabon's total cases now exceed
#   Croacia's, so the two swap places in the (descending, by total cases)
#   sorted table: row 87 becomes Gabon, row 88 becomes Croacia.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Datos actualizados a ..." timestamp -------------------------
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 23:35"

# --- Estados Unidos (row 4) ------------------------------------------------
$ws.Cells.Item(4, 2).Value = 1741808   # Casos totales
$ws.Cells.Item(4, 3).Value = 16533     # Nuevos casos
$ws.Cells.Item(4, 4).Value = 485300    # Casos activos
$ws.Cells.Item(4, 5).Value = 1154679   # Recuperados
$ws.Cells.Item(4, 7).Value = 1257      # Muertes hoy
$ws.Cells.Item(4, 8).Value = 101829    # Muertes

# --- Irlanda (row 34) -------------------------------------------------------
$ws.Cells.Item(34, 4).Value = 22089    # Casos activos
$ws.Cells.Item(34, 5).Value = 1083     # Recuperados

# --- Gabon overtakes Croacia: row 87 now holds Gabon's (updated) data ------
$ws.Cells.Item(87, 1).Value = "Gabon"
$ws.Cells.Item(87, 2).Value = 2319     # Casos totales
$ws.Cells.Item(87, 3).Value = 81       # Nuevos casos
$ws.Cells.Item(87, 4).Value = 631      # Casos activos
$ws.Cells.Item(87, 5).Value = 1674     # Recuperados
$ws.Cells.Item(87, 8).Value = 14       # Muertes

# --- ... and row 88 now holds Croacia's (unchanged) data -------------------
$ws.Cells.Item(88, 1).Value = "Croacia"
$ws.Cells.Item(88, 2).Value = 2244     # Casos totales
$ws.Cells.Item(88, 4).Value = 2047     # Casos activos
$ws.Cells.Item(88, 5).Value = 96       # Recuperados
$ws.Cells.Item(88, 8).Value = 101      # Muertes
